# Insert a new weekly record at the top of the "Haba" / Feria Lagunitas de
# Puerto Montt data block (row 138), pushing all the existing rows down by
# one. This mirrors the "Fruta / hortaliza, semanal" weekly refresh: a new
# price observation is prepended and the dimension grows from R170 to R171.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 138:170 down to 139:171 and leave a blank row 138.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Cells.Item(138, 1).Value  = 4
$ws.Cells.Item(138, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(138, 3).Value  = "Los Lagos"
$ws.Cells.Item(138, 4).Value  = 45218
$ws.Cells.Item(138, 5).Value  = 10
$ws.Cells.Item(138, 6).Value  = 100112026
$ws.Cells.Item(138, 7).Value  = "Haba"
$ws.Cells.Item(138, 8).Value  = "Sin especificar"
$ws.Cells.Item(138, 9).Value  = "Primera"
$ws.Cells.Item(138, 10).Value = 60
$ws.Cells.Item(138, 11).Value = 16000
$ws.Cells.Item(138, 12).Value = 16000
$ws.Cells.Item(138, 13).Value = 16000
$ws.Cells.Item(138, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(138, 15).Value = "Región Metropolitana"
$ws.Cells.Item(138, 16).Value = 640
$ws.Cells.Item(138, 17).Value = 25
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date style used by the rest of
# column D (it is inherited from the insert, but set it explicitly too).
$ws.Cells.Item(138, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
